# Weekly Fruta/Hortaliza update: insert a new price record as row 97
# (Haba, Vega Modelo de Temuco) and push the existing rows 97-100 down
# to rows 98-101, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 97; this shifts the old
# rows 97-100 down to 98-101, growing the sheet from 100 to 101 rows.
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new record.
$ws.Cells.Item(97, 1).Value = 10
$ws.Cells.Item(97, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(97, 3).Value = "La Araucanía"
$ws.Cells.Item(97, 4).Value = 45239
$ws.Cells.Item(97, 5).Value = 9
$ws.Cells.Item(97, 6).Value = 100112026
$ws.Cells.Item(97, 7).Value = "Haba"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 220
$ws.Cells.Item(97, 11).Value = 10000
$ws.Cells.Item(97, 12).Value = 12000
$ws.Cells.Item(97, 13).Value = 11091
$ws.Cells.Item(97, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(97, 15).Value = "Región Metropolitana"
$ws.Cells.Item(97, 16).Value = 444
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Keep the date cell formatted the same way as the other date cells in
# column D (style index 2, yyyy-mm-dd hh:mm:ss number format).
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
